$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 118 ---
$ws.Range("A118").Value = 117
$ws.Range("B118").Value = "paraguay"
$ws.Range("C118").Value = "primera-division"
$ws.Range("D118").NumberFormat = "@"
$ws.Range("D118").Value = "2023"
$ws.Range("D118").ClearFormats()
$ws.Range("E118").Value = 45242.95833333334
$ws.Range("F118").Value = "Cerro Porteno"
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = "Tacuary"
$ws.Range("I118").Value = 1
$ws.Range("J118").Value = 1.3
$ws.Range("K118").Value = "08/11/2023 14:42"
$ws.Range("L118").Value = 1.33
$ws.Range("M118").Value = "12/11/2023 22:59"
$ws.Range("N118").Value = 5.79
$ws.Range("O118").Value = "08/11/2023 14:42"
$ws.Range("P118").Value = 5.58
$ws.Range("Q118").Value = "12/11/2023 22:59"
$ws.Range("R118").Value = 9.390000000000001
$ws.Range("S118").Value = "08/11/2023 14:42"
$ws.Range("T118").Value = 9.119999999999999
$ws.Range("U118").Value = "12/11/2023 22:59"
$ws.Range("V118").Value = "https://www.betexplorer.com/football/paraguay/primera-division/cerro-porteno-tacuary/08MNZtck/"

# --- Row 119 ---
$ws.Range("A119").Value = 118
$ws.Range("B119").Value = "paraguay"
$ws.Range("C119").Value = "primera-division"
$ws.Range("D119").NumberFormat = "@"
$ws.Range("D119").Value = "2023"
$ws.Range("D119").ClearFormats()
$ws.Range("E119").Value = 45242.95833333334
$ws.Range("F119").Value = "Sp. Luqueno"
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = "Libertad Asuncion"
$ws.Range("I119").Value = 1
$ws.Range("J119").Value = 4.03
$ws.Range("K119").Value = "08/11/2023 14:42"
$ws.Range("L119").Value = 4.08
$ws.Range("M119").Value = "12/11/2023 22:58"
$ws.Range("N119").Value = 3.64
$ws.Range("O119").Value = "08/11/2023 14:42"
$ws.Range("P119").Value = 3.38
$ws.Range("Q119").Value = "12/11/2023 22:58"
$ws.Range("R119").Value = 1.92
$ws.Range("S119").Value = "08/11/2023 14:42"
$ws.Range("T119").Value = 2.03
$ws.Range("U119").Value = "12/11/2023 22:58"
$ws.Range("V119").Value = "https://www.betexplorer.com/football/paraguay/primera-division/sp-luqueno-libertad-asuncion/rXBSY0Ce/"

# --- Copy cell formatting (style) from row 117 template so new rows match
#     the rest of the table: column A uses the bold/bordered "index" style,
#     column E uses the datetime display style.
$ws.Range("A117").Copy()
$ws.Range("A118:A119").PasteSpecial(-4122)

$ws.Range("E117").Copy()
$ws.Range("E118:E119").PasteSpecial(-4122)

$excel.CutCopyMode = 0
